$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.77148166666667
$ws.Range("H2").Value = 107.314445
$ws.Range("I2").Value = 0.1058641704420874
$ws.Range("J2").Value = 0.1112463097643854
$ws.Range("M2").Value = 2.724367666666667
$ws.Range("N2").Value = 8.173103000000001
$ws.Range("O2").Value = 0.0714746893937865
$ws.Range("P2").Value = 0.08141055085372029
$ws.Range("Q2").Value = 97.45466804142613
$ws.Range("R2").Value = 877.0920123728351
$ws.Range("S2").Value = 0.007566608700279067
$ws.Range("T2").Value = 0.009056623358362215
$ws.Range("G3").Value = 35.77148166666667
$ws.Range("H3").Value = 107.314445
$ws.Range("I3").Value = 0.1058641704420874
$ws.Range("J3").Value = 0.1112463097643854
$ws.Range("O3").Value = 0.3977336888558922
$ws.Range("P3").Value = 0.4530235664886322
$ws.Range("Q3").Value = 542.3039252789316
$ws.Range("R3").Value = 4880.735327510385
$ws.Range("S3").Value = 0.04210574702760031
$ws.Range("T3").Value = 0.05039720000816101
$ws.Range("G4").Value = 35.77148166666667
$ws.Range("H4").Value = 107.314445
$ws.Range("I4").Value = 0.1058641704420874
$ws.Range("J4").Value = 0.1112463097643854
$ws.Range("M4").Value = 3.339883333333333
$ws.Range("N4").Value = 10.01965
$ws.Range("O4").Value = 0.0876229470721772
$ws.Range("P4").Value = 0.09980361508492896
$ws.Range("Q4").Value = 119.4725754271389
$ws.Range("R4").Value = 1075.25317884425
$ws.Range("S4").Value = 0.009276130603486967
$ws.Range("T4").Value = 0.01110278387934349
$ws.Range("G5").Value = 35.77148166666667
$ws.Range("H5").Value = 107.314445
$ws.Range("I5").Value = 0.1058641704420874
$ws.Range("J5").Value = 0.1112463097643854
$ws.Range("M5").Value = 13.955954
$ws.Range("N5").Value = 27.911908
$ws.Range("O5").Value = 0.366139082308386
$ws.Range("P5").Value = 0.2780246138655492
$ws.Range("Q5").Value = 499.2251526518434
$ws.Range("R5").Value = 2995.35091591106
$ws.Range("S5").Value = 0.03876101021500443
$ws.Range("T5").Value = 0.03092921231621052
$ws.Range("G6").Value = 35.77148166666667
$ws.Range("H6").Value = 107.314445
$ws.Range("I6").Value = 0.1058641704420874
$ws.Range("J6").Value = 0.1112463097643854
$ws.Range("M6").Value = 2.936101333333333
$ws.Range("N6").Value = 8.808304
$ws.Range("O6").Value = 0.07702959236975811
$ws.Range("P6").Value = 0.08773765370716943
$ws.Range("Q6").Value = 105.0286950168089
$ws.Range("R6").Value = 945.25825515128
$ws.Range("S6").Value = 0.008154673895716584
$ws.Range("T6").Value = 0.009760490202308146
$ws.Range("I7").Value = 0.4187622210170216
$ws.Range("J7").Value = 0.4400521117044616
$ws.Range("M7").Value = 2.724367666666667
$ws.Range("N7").Value = 8.173103000000001
$ws.Range("O7").Value = 0.0714746893937865
$ws.Range("P7").Value = 0.08141055085372029
$ws.Range("Q7").Value = 385.4971239757583
$ws.Range("R7").Value = 3469.474115781824
$ws.Range("S7").Value = 0.02993089967704379
$ws.Range("T7").Value = 0.03582488481820308
$ws.Range("I8").Value = 0.4187622210170216
$ws.Range("J8").Value = 0.4400521117044616
$ws.Range("O8").Value = 0.3977336888558922
$ws.Range("P8").Value = 0.4530235664886322
$ws.Range("Q8").Value = 2145.167673516949
$ws.Range("S8").Value = 0.1665558429185864
$ws.Range("T8").Value = 0.1993539770852092
$ws.Range("I9").Value = 0.4187622210170216
$ws.Range("J9").Value = 0.4400521117044616
$ws.Range("M9").Value = 3.339883333333333
$ws.Range("N9").Value = 10.01965
$ws.Range("O9").Value = 0.0876229470721772
$ws.Range("P9").Value = 0.09980361508492896
$ws.Range("Q9").Value = 472.5923872785777
$ws.Range("R9").Value = 4253.331485507199
$ws.Range("S9").Value = 0.03669317992800186
$ws.Range("T9").Value = 0.04391879157386226
$ws.Range("I10").Value = 0.4187622210170216
$ws.Range("J10").Value = 0.4400521117044616
$ws.Range("M10").Value = 13.955954
$ws.Range("N10").Value = 27.911908
$ws.Range("O10").Value = 0.366139082308386
$ws.Range("P10").Value = 0.2780246138655492
$ws.Range("Q10").Value = 1974.762876231211
$ws.Range("R10").Value = 11848.57725738726
$ws.Range("S10").Value = 0.1533252153085938
$ws.Range("T10").Value = 0.1223453184373525
$ws.Range("I11").Value = 0.4187622210170216
$ws.Range("J11").Value = 0.4400521117044616
$ws.Range("M11").Value = 2.936101333333333
$ws.Range("N11").Value = 8.808304
$ws.Range("O11").Value = 0.07702959236975811
$ws.Range("P11").Value = 0.08773765370716943
$ws.Range("Q11").Value = 415.4573677958258
$ws.Range("R11").Value = 3739.116310162432
$ws.Range("S11").Value = 0.03225708318479573
$ws.Range("T11").Value = 0.03860913978983469
$ws.Range("G12").Value = 52.33127733333333
$ws.Range("H12").Value = 156.993832
$ws.Range("I12").Value = 0.1548721776383825
$ws.Range("J12").Value = 0.1627458863135329
$ws.Range("M12").Value = 2.724367666666667
$ws.Range("N12").Value = 8.173103000000001
$ws.Range("O12").Value = 0.0714746893937865
$ws.Range("P12").Value = 0.08141055085372029
$ws.Range("Q12").Value = 142.5696399222996
$ws.Range("R12").Value = 1283.126759300696
$ws.Range("S12").Value = 0.01106944079244272
$ws.Range("T12").Value = 0.01324923225396165
$ws.Range("G13").Value = 52.33127733333333
$ws.Range("H13").Value = 156.993832
$ws.Range("I13").Value = 0.1548721776383825
$ws.Range("J13").Value = 0.1627458863135329
$ws.Range("O13").Value = 0.3977336888558922
$ws.Range("P13").Value = 0.4530235664886322
$ws.Range("Q13").Value = 793.3542528983972
$ws.Range("R13").Value = 7140.188276085575
$ws.Range("S13").Value = 0.06159788251325889
$ws.Range("T13").Value = 0.07372772184911014
$ws.Range("G14").Value = 52.33127733333333
$ws.Range("H14").Value = 156.993832
$ws.Range("I14").Value = 0.1548721776383825
$ws.Range("J14").Value = 0.1627458863135329
$ws.Range("M14").Value = 3.339883333333333
$ws.Range("N14").Value = 10.01965
$ws.Range("O14").Value = 0.0876229470721772
$ws.Range("P14").Value = 0.09980361508492896
$ws.Range("Q14").Value = 174.7803609776444
$ws.Range("R14").Value = 1573.0232487988
$ws.Range("S14").Value = 0.01357035662416082
$ws.Range("T14").Value = 0.01624262779429144
$ws.Range("G15").Value = 52.33127733333333
$ws.Range("H15").Value = 156.993832
$ws.Range("I15").Value = 0.1548721776383825
$ws.Range("J15").Value = 0.1627458863135329
$ws.Range("M15").Value = 13.955954
$ws.Range("N15").Value = 27.911908
$ws.Range("O15").Value = 0.366139082308386
$ws.Range("P15").Value = 0.2780246138655492
$ws.Range("Q15").Value = 730.3328992252426
$ws.Range("R15").Value = 4381.997395351456
$ws.Range("S15").Value = 0.05670475699561871
$ws.Range("T15").Value = 0.04524736220052655
$ws.Range("G16").Value = 52.33127733333333
$ws.Range("H16").Value = 156.993832
$ws.Range("I16").Value = 0.1548721776383825
$ws.Range("J16").Value = 0.1627458863135329
$ws.Range("M16").Value = 2.936101333333333
$ws.Range("N16").Value = 8.808304
$ws.Range("O16").Value = 0.07702959236975811
$ws.Range("P16").Value = 0.08773765370716943
$ws.Range("Q16").Value = 153.6499331534365
$ws.Range("R16").Value = 1382.849398380928
$ws.Range("S16").Value = 0.01192974071290137
$ws.Range("T16").Value = 0.01427894221564311
$ws.Range("G17").Value = 49.043167
$ws.Range("H17").Value = 98.08633399999999
$ws.Range("I17").Value = 0.1451411557029742
$ws.Range("J17").Value = 0.1016800925151965
$ws.Range("M17").Value = 2.724367666666667
$ws.Range("N17").Value = 8.173103000000001
$ws.Range("O17").Value = 0.0714746893937865
$ws.Range("P17").Value = 0.08141055085372029
$ws.Range("Q17").Value = 133.6116184457337
$ws.Range("R17").Value = 801.669710674402
$ws.Range("S17").Value = 0.01037391902212528
$ws.Range("T17").Value = 0.008277832342519385
$ws.Range("G18").Value = 49.043167
$ws.Range("H18").Value = 98.08633399999999
$ws.Range("I18").Value = 0.1451411557029742
$ws.Range("J18").Value = 0.1016800925151965
$ws.Range("O18").Value = 0.3977336888558922
$ws.Range("P18").Value = 0.4530235664886322
$ws.Range("Q18").Value = 743.5057406915769
$ws.Range("R18").Value = 4461.034444149461
$ws.Range("S18").Value = 0.05772752726255134
$ws.Range("T18").Value = 0.04606347815212838
$ws.Range("G19").Value = 49.043167
$ws.Range("H19").Value = 98.08633399999999
$ws.Range("I19").Value = 0.1451411557029742
$ws.Range("J19").Value = 0.1016800925151965
$ws.Range("M19").Value = 3.339883333333333
$ws.Range("N19").Value = 10.01965
$ws.Range("O19").Value = 0.0876229470721772
$ws.Range("P19").Value = 0.09980361508492896
$ws.Range("Q19").Value = 163.7984560771833
$ws.Range("R19").Value = 982.7907364630998
$ws.Range("S19").Value = 0.01271769580415634
$ws.Range("T19").Value = 0.01014804081518664
$ws.Range("G20").Value = 49.043167
$ws.Range("H20").Value = 98.08633399999999
$ws.Range("I20").Value = 0.1451411557029742
$ws.Range("J20").Value = 0.1016800925151965
$ws.Range("M20").Value = 13.955954
$ws.Range("N20").Value = 27.911908
$ws.Range("O20").Value = 0.366139082308386
$ws.Range("P20").Value = 0.2780246138655492
$ws.Range("Q20").Value = 684.4441826663179
$ws.Range("R20").Value = 2737.776730665272
$ws.Range("S20").Value = 0.05314184955426553
$ws.Range("T20").Value = 0.02826956845935082
$ws.Range("G21").Value = 49.043167
$ws.Range("H21").Value = 98.08633399999999
$ws.Range("I21").Value = 0.1451411557029742
$ws.Range("J21").Value = 0.1016800925151965
$ws.Range("M21").Value = 2.936101333333333
$ws.Range("N21").Value = 8.808304
$ws.Range("O21").Value = 0.07702959236975811
$ws.Range("P21").Value = 0.08773765370716943
$ws.Range("Q21").Value = 143.9957080195893
$ws.Range("R21").Value = 863.9742481175359
$ws.Range("S21").Value = 0.01118016405987569
$ws.Range("T21").Value = 0.008921172746011258
$ws.Range("G22").Value = 59.25420133333333
$ws.Range("H22").Value = 177.762604
$ws.Range("I22").Value = 0.1753602751995342
$ws.Range("J22").Value = 0.1842755997024237
$ws.Range("M22").Value = 2.724367666666667
$ws.Range("N22").Value = 8.173103000000001
$ws.Range("O22").Value = 0.0714746893937865
$ws.Range("P22").Value = 0.08141055085372029
$ws.Range("Q22").Value = 161.4302302266903
$ws.Range("R22").Value = 1452.872072040212
$ws.Range("S22").Value = 0.01253382120189563
$ws.Range("T22").Value = 0.01500197808067397
$ws.Range("G23").Value = 59.25420133333333
$ws.Range("H23").Value = 177.762604
$ws.Range("I23").Value = 0.1753602751995342
$ws.Range("J23").Value = 0.1842755997024237
$ws.Range("O23").Value = 0.3977336888558922
$ws.Range("P23").Value = 0.4530235664886322
$ws.Range("Q23").Value = 898.3073799338412
$ws.Range("R23").Value = 8084.766419404572
$ws.Range("S23").Value = 0.06974668913389517
$ws.Range("T23").Value = 0.08348118939402353
$ws.Range("G24").Value = 59.25420133333333
$ws.Range("H24").Value = 177.762604
$ws.Range("I24").Value = 0.1753602751995342
$ws.Range("J24").Value = 0.1842755997024237
$ws.Range("M24").Value = 3.339883333333333
$ws.Range("N24").Value = 10.01965
$ws.Range("O24").Value = 0.0876229470721772
$ws.Range("P24").Value = 0.09980361508492896
$ws.Range("Q24").Value = 197.9021194631777
$ws.Range("R24").Value = 1781.1190751686
$ws.Range("S24").Value = 0.01536558411237122
$ws.Range("T24").Value = 0.01839137102224515
$ws.Range("G25").Value = 59.25420133333333
$ws.Range("H25").Value = 177.762604
$ws.Range("I25").Value = 0.1753602751995342
$ws.Range("J25").Value = 0.1842755997024237
$ws.Range("M25").Value = 13.955954
$ws.Range("N25").Value = 27.911908
$ws.Range("O25").Value = 0.366139082308386
$ws.Range("P25").Value = 0.2780246138655492
$ws.Range("Q25").Value = 826.9489081147387
$ws.Range("R25").Value = 4961.693448688432
$ws.Range("S25").Value = 0.06420625023490348
$ws.Range("T25").Value = 0.05123315245210888
$ws.Range("G26").Value = 59.25420133333333
$ws.Range("H26").Value = 177.762604
$ws.Range("I26").Value = 0.1753602751995342
$ws.Range("J26").Value = 0.1842755997024237
$ws.Range("M26").Value = 2.936101333333333
$ws.Range("N26").Value = 8.808304
$ws.Range("O26").Value = 0.07702959236975811
$ws.Range("P26").Value = 0.08773765370716943
$ws.Range("Q26").Value = 173.9763395404018
$ws.Range("R26").Value = 1565.787055863616
$ws.Range("S26").Value = 0.01350793051646872
$ws.Range("T26").Value = 0.01616790875337223
